$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Copyright 2017 John Brzezicki*") {
        $p.Range.Delete()
        break
    }
}
